$wb = $excel.ActiveWorkbook

# --- Sheet "discard": remove the "amused" entry (row 1 shifts everything up) ---
$ws5 = $wb.Worksheets.Item("discard")
$ws5.Rows.Item(1).Delete()
$ws5.Rows.Item(1).EntireRow.Select() | Out-Null

# --- Sheet "samples_retained": add EYASE dataset row + trailing helper row ---
$ws1 = $wb.Worksheets.Item("samples_retained")

$ws1.Cells.Item(18,1).Value = "EYASE"
$ws1.Cells.Item(18,2).Value = "acted"
$ws1.Cells.Item(18,3).Value = 132
$ws1.Cells.Item(18,4).Value = 297
$ws1.Cells.Item(18,5).Value = 150
$ws1.Cells.Item(18,6).Value = "Arabic"
$ws1.Cells.Item(18,7).Formula = "=IF(OR(ISBLANK(C18), ISBLANK(D18),ISBLANK(E18)), """", SUM(C18:E18))"
$ws1.Cells.Item(18,8).Value = "Egyptian Arabic from a TV drama"

$ws1.Cells.Item(19,7).Formula = "=IF(OR(ISBLANK(C19), ISBLANK(D19),ISBLANK(E19)), """", SUM(C19:E19))"

$ws1.Application.Goto($ws1.Range("A19"))
